# Apply updated cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.108.55"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "'1.743.66"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'243.30"
$ws.Range("E5").Value = "  +5.30%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.5320"
$ws.Range("E7").Value = "  +2.95%  "

$ws.Range("D8").Value = "'0.2801"
$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("D9").Value = "'0.06184"
$ws.Range("E9").Value = "  +1.43%  "

$ws.Range("D10").Value = "'1.745.05"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").Value = "'0.07197"
$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("D12").Value = "'15.37"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").Value = "'0.6567"
$ws.Range("E13").Value = "  +2.82%  "

$ws.Range("D14").Value = "'4.643"
$ws.Range("E14").Value = "  +3.13%  "

$ws.Range("D15").Value = "'77.96"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "'26.130.67"
$ws.Range("E18").Value = "  +1.05%  "

$ws.Range("D19").Value = "'11.90"
$ws.Range("E19").Value = "  +4.15%  "

$ws.Range("D20").Value = "'0.000006790"

$ws.Range("D21").Value = "'1.965.62"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").Value = "'4.453"
$ws.Range("E22").Value = "  +7.95%  "

$ws.Range("D23").Value = "'8.746"
$ws.Range("E23").Value = "  +1.60%  "

$ws.Range("D24").Value = "'5.266"
$ws.Range("E24").Value = "  +2.62%  "

$ws.Range("D25").Value = "'140.95"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").Value = "'15.33"
$ws.Range("E27").Value = "  +1.92%  "

$ws.Range("D28").Value = "'1.792"
$ws.Range("E28").Value = "  -0.80%  "

$ws.Range("D29").Value = "'105.89"
$ws.Range("E29").Value = "  +3.78%  "

$ws.Range("D30").Value = "'0.08463"
$ws.Range("E30").Value = "  +2.77%  "

$ws.Range("D31").Value = "'3.880"
$ws.Range("E31").Value = "  +6.18%  "

$ws.Range("E32").Value = "  +7.46%  "

$ws.Range("D33").Value = "'0.04623"
$ws.Range("E33").Value = "  +3.15%  "

$ws.Range("D34").Value = "'2.667"
$ws.Range("E34").Value = "  +2.26%  "

$ws.Range("D35").Value = "'0.9996"
$ws.Range("E35").Value = "  +2.11%  "

$ws.Range("D36").Value = "'0.6304"
$ws.Range("E36").Value = "  +2.83%  "

$ws.Range("D37").Value = "'2.707"
$ws.Range("E37").Value = "  +2.26%  "

$ws.Range("D38").Value = "'0.01620"
$ws.Range("E38").Value = "  +2.25%  "

$ws.Range("D39").Value = "'1.947"
$ws.Range("E39").Value = "  +1.51%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").Value = "'99.96"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("D42").Value = "'0.3921"
$ws.Range("E42").Value = "  +2.49%  "

$ws.Range("D43").Value = "'0.7533"
$ws.Range("E43").Value = "  +4.48%  "

$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("D45").Value = "'0.1152"
$ws.Range("E45").Value = "  +3.03%  "

$ws.Range("D46").Value = "'6.316"
$ws.Range("E46").Value = "  +0.85%  "

$ws.Range("D47").Value = "'0.05332"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("D48").Value = "'55.10"
$ws.Range("E48").Value = "  +3.81%  "

$ws.Range("D49").Value = "'30.97"
$ws.Range("E49").Value = "  +4.10%  "

$ws.Range("D50").Value = "'0.3489"
$ws.Range("E50").Value = "  +3.34%  "

$ws.Range("D51").Value = "'7.583"
$ws.Range("E51").Value = "  -0.57%  "
